# Automatische test-sync: 2025-06-29 15:00:50
# Adds Testmail #9 ("Wat zijn jullie voorwaarden?") as a new row to the
# "Logs" sheet and updates the category summary on the "Dashboard" sheet
# to reflect the new Productinformatie count / row order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new testmail entry as row 24
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 24

$logs.Cells.Item($newRow, 1).Value = "Wat zijn jullie voorwaarden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #9: Wat zijn jullie voorwaarden?"
$logs.Cells.Item($newRow, 4).Value = "Productinformatie"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nDank voor uw interesse in onze diensten. Voor informatie over onze voorwaarden kunt u terecht op onze website onder de sectie 'Algemene Voorwaarden'. Mocht u specifieke vragen hebben, dan helpen wij u graag verder. Aarzel niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-29 15:00:07"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Nee"
$logs.Cells.Item($newRow, 9).Value = "Ja"

# Restore the row height back to the sheet default - setting a wrapped,
# multi-line value otherwise leaves an explicit autofit height behind.
$logs.Rows.Item($newRow).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend conditional formatting ranges (D/G/H/I) from row 23 to row 24
#    without disturbing the existing rules themselves.
# ---------------------------------------------------------------------
$logs.Range("D2:D23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D24"))
$logs.Range("G2:G23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G24"))
$logs.Range("H2:H23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H24"))
$logs.Range("I2:I23").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I24"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: swap rows 3/4 (category order) and bump the
#    Productinformatie count from 5 to 6
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(3, 1).Value = "Productinformatie"
$dash.Cells.Item(3, 2).Value = 6
$dash.Cells.Item(4, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(4, 2).Value = 5

$wb.Save()
